$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 11.14494766666667
$ws.Range("H2").Value = 33.434843
$ws.Range("I2").Value = 0.1279818847384872
$ws.Range("J2").Value = 0.1279818847384872
$ws.Range("M2").Value = 0.01392333333333333
$ws.Range("N2").Value = 0.04177
$ws.Range("O2").Value = 0.008343913876905598
$ws.Range("P2").Value = 0.008343913876905598
$ws.Range("Q2").Value = 0.1551748213455556
$ws.Range("R2").Value = 1.39657339211
$ws.Range("S2").Value = 0.001067869824061997
$ws.Range("T2").Value = 0.001067869824061997
$ws.Range("G3").Value = 11.14494766666667
$ws.Range("H3").Value = 33.434843
$ws.Range("I3").Value = 0.1279818847384872
$ws.Range("J3").Value = 0.1279818847384872
$ws.Range("O3").Value = 0.09015262350870268
$ws.Range("P3").Value = 0.09015262350870269
$ws.Range("Q3").Value = 1.676601347182667
$ws.Range("R3").Value = 15.089412124644
$ws.Range("S3").Value = 0.01153790267076302
$ws.Range("T3").Value = 0.01153790267076302
$ws.Range("G4").Value = 11.14494766666667
$ws.Range("H4").Value = 33.434843
$ws.Range("I4").Value = 0.1279818847384872
$ws.Range("J4").Value = 0.1279818847384872
$ws.Range("M4").Value = 1.504322
$ws.Range("N4").Value = 4.512966
$ws.Range("O4").Value = 0.9015034626143917
$ws.Range("P4").Value = 0.9015034626143917
$ws.Range("Q4").Value = 16.76558996381533
$ws.Range("R4").Value = 150.890309674338
$ws.Range("S4").Value = 0.1153761122436622
$ws.Range("T4").Value = 0.1153761122436622
$ws.Range("I5").Value = 0.5307607770439682
$ws.Range("J5").Value = 0.5307607770439681
$ws.Range("M5").Value = 0.01392333333333333
$ws.Range("N5").Value = 0.04177
$ws.Range("O5").Value = 0.008343913876905598
$ws.Range("P5").Value = 0.008343913876905598
$ws.Range("Q5").Value = 0.6435341136233333
$ws.Range("R5").Value = 5.79180702261
$ws.Range("S5").Value = 0.004428622212894365
$ws.Range("T5").Value = 0.004428622212894364
$ws.Range("I6").Value = 0.5307607770439682
$ws.Range("J6").Value = 0.5307607770439681
$ws.Range("O6").Value = 0.09015262350870268
$ws.Range("P6").Value = 0.09015262350870269
$ws.Range("S6").Value = 0.04784947650603134
$ws.Range("T6").Value = 0.04784947650603134
$ws.Range("I7").Value = 0.5307607770439682
$ws.Range("J7").Value = 0.5307607770439681
$ws.Range("M7").Value = 1.504322
$ws.Range("N7").Value = 4.512966
$ws.Range("O7").Value = 0.9015034626143917
$ws.Range("P7").Value = 0.9015034626143917
$ws.Range("Q7").Value = 69.529508609582
$ws.Range("R7").Value = 625.7655774862379
$ws.Range("S7").Value = 0.4784826783250425
$ws.Range("T7").Value = 0.4784826783250424
$ws.Range("G8").Value = 29.71744933333333
$ws.Range("H8").Value = 89.152348
$ws.Range("I8").Value = 0.3412573382175446
$ws.Range("J8").Value = 0.3412573382175446
$ws.Range("M8").Value = 0.01392333333333333
$ws.Range("N8").Value = 0.04177
$ws.Range("O8").Value = 0.008343913876905598
$ws.Range("P8").Value = 0.008343913876905598
$ws.Range("Q8").Value = 0.4137659528844445
$ws.Range("R8").Value = 3.72389357596
$ws.Range("S8").Value = 0.002847421839949237
$ws.Range("T8").Value = 0.002847421839949237
$ws.Range("G9").Value = 29.71744933333333
$ws.Range("H9").Value = 89.152348
$ws.Range("I9").Value = 0.3412573382175446
$ws.Range("J9").Value = 0.3412573382175446
$ws.Range("O9").Value = 0.09015262350870268
$ws.Range("P9").Value = 0.09015262350870269
$ws.Range("Q9").Value = 4.470574207909333
$ws.Range("R9").Value = 40.235167871184
$ws.Range("S9").Value = 0.03076524433190831
$ws.Range("T9").Value = 0.03076524433190832
$ws.Range("G10").Value = 29.71744933333333
$ws.Range("H10").Value = 89.152348
$ws.Range("I10").Value = 0.3412573382175446
$ws.Range("J10").Value = 0.3412573382175446
$ws.Range("M10").Value = 1.504322
$ws.Range("N10").Value = 4.512966
$ws.Range("O10").Value = 0.9015034626143917
$ws.Range("P10").Value = 0.9015034626143917
$ws.Range("Q10").Value = 44.70461281601867
$ws.Range("R10").Value = 402.341515344168
$ws.Range("S10").Value = 0.3076446720456871
$ws.Range("T10").Value = 0.3076446720456871
